$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the timestamp in A73 (2024-06-05 row): 07:25:12.999996 -> 07:00:00 ---
$ws.Range("A73").Value = 45448.2916666667

# --- Append a new row 74 for 2024-06-06 ---
$ws.Range("A74").Value = 45449.5011111111
$ws.Range("B74").Value = 1500
$ws.Range("C74").Value = 2.97000002861023
$ws.Range("D74").Value = 2.97000002861023
$ws.Range("E74").Value = 2.97000002861023
$ws.Range("F74").Value = 2.97000002861023

# G74 / H74 hold numeric-looking text in this workbook (shared string that
# equals the close price, and the ticker). Force text typing via a
# temporary "@" number format, then clear the format stamp so the cell
# keeps its default style like its neighbours.
$ws.Range("G74").NumberFormat = "@"
$ws.Range("G74").Value = "2.97000002861023"
$ws.Range("G74").ClearFormats()

$ws.Range("H74").NumberFormat = "@"
$ws.Range("H74").Value = "ESPE.MI"
$ws.Range("H74").ClearFormats()

# A74 should carry the same date/time style as the rest of column A.
$ws.Range("A73").Copy()
$ws.Range("A74").PasteSpecial(-4122) | Out-Null
$ws.Range("A74").Value = 45449.5011111111
